$d = $word.ActiveDocument

# Locate the paragraph that currently reads (Hebrew RTL text, anchor is ASCII):
# "Waitpid() - <hebrew text ending in> ...lesimush."
# by finding the unique anchor text "Waitpid(" and expanding to the whole paragraph.
$rng = $d.Content
$rng.Find.Execute("Waitpid(", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Expand(4) | Out-Null

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:rtl/></w:rPr></w:pPr>
<w:proofErr w:type="spellStart"/>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>Waitpid</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t>(</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t>)</w:t></w:r>
<w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:rPr><w:rtl/></w:rPr><w:t>&#8211;</w:t></w:r>
<w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve">&#1495;&#1505;&#1512;&#1493;&#1503;: </w:t></w:r>
<w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>&#1512;&#1511; &#1506;&#1489;&#1493;&#1512; &#1514;&#1492;&#1500;&#1497;&#1498; &#1513;&#1492;&#1493;&#1488; &#1489;&#1503; &#1513;&#1500; &#1492;&#1514;&#1492;&#1500;&#1497;&#1498; &#1492;&#1504;&#1493;&#1499;&#1495;&#1497;, &#1497;&#1514;&#1512;&#1493;&#1503;: &#1492;&#1493;&#1488; &#1489;&#1496;&#1493;&#1495; &#1500;&#1513;&#1497;&#1502;&#1493;&#1513;.</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve"> (&#1499;&#1504;&#1512;&#1488;&#1492; &#1513;&#1510;&#1512;&#1497;&#1498; &#1500;&#1492;&#1513;&#1514;&#1502;&#1513; &#1489;&#1488;&#1493;&#1508;&#1510;&#1497;&#1492; &#1513;&#1500; </w:t></w:r>
<w:r><w:rPr><w:rFonts w:hint="cs"/></w:rPr><w:t>WNOHANG</w:t></w:r>
<w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>)</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:rFonts w:hint="cs"/></w:rPr></w:pPr>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$rng.InsertXML($xml) | Out-Null
